# Apply the "Add new JS files and updated a few" edit:
# Adds two new API-documentation blocks to sheet1 of the workbook:
#   1) app_group_user_apply_for_member_MS   (rows 35-39, after the existing row 34 blank separator)
#   2) sys_sec_group_owner_approve_user_membership_MS (rows 40-44)
#
# The pre-existing blank separator rows 34/35 (B/C cells styled "Header3")
# are extended: row 34 grows to include D/E wrap-style cells, and row 35
# becomes the first new section title row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- widen column A slightly (66 -> 71), matching the diff's <col> change ---
$ws.Columns.Item(1).ColumnWidth = 71

# --- row 34: existing blank separator row gains wrap-style D/E cells ---
$ws.Range("D34").Style = "Header4"
$ws.Range("E34").Style = "Header4"

# ============================================================
# Section: app_group_user_apply_for_member_MS  (rows 35-39)
# ============================================================

# Row 35: section title + short description (mirrors existing sections)
$ws.Range("A35").Value = "app_group_user_apply_for_member_MS"
$ws.Range("A35").Style = "Header1"
$ws.Range("B35").Style = "Header3"
$ws.Range("C35").Value = "用户申请入群"
$ws.Range("C35").Style = "Header2"
$ws.Range("D35").Style = "Header4"
$ws.Range("E35").Style = "Header4"
$ws.Rows.Item(35).RowHeight = 21

# Row 36: "Parameters" label row
$ws.Range("A36").Value = "Parameters"
$ws.Range("A36").Style = "Bold12"
$ws.Range("B36").Style = "Header3"
$ws.Range("C36").Style = "Header3"
$ws.Range("D36").Style = "Header4"
$ws.Range("E36").Style = "Header4"

# Row 37: Parameter 1 - @usertToken / 用户·Token
$ws.Range("A37").Value = "    Parameter 1"
$ws.Range("A37").Style = "Header3"
$ws.Range("B37").Value = "@usertToken"
$ws.Range("B37").Style = "Header3"
$ws.Range("C37").Value = "用户·Token"
$ws.Range("C37").Style = "Header3"
$ws.Range("D37").Style = "Header4"
$ws.Range("E37").Style = "Header4"

# Row 38: Parameter 2 - @groupName / 群名字 (+ temporary notes, row taller/wraps)
$ws.Range("A38").Value = "    Parameter 2"
$ws.Range("A38").Style = "Header3"
$ws.Range("B38").Value = "@groupName"
$ws.Range("C38").Value = "群名字"
$ws.Range("D38").Value = "界面上显示可以加入的群·名 （临时）： Fight For Trump"
$ws.Range("D38").Style = "Wrap"
$ws.Rows.Item(38).RowHeight = 45

# Row 39: trailing blank separator row for this section
$ws.Range("A39").Style = "Header3"
$ws.Range("B39").Style = "Header3"
$ws.Range("C39").Style = "Header3"
$ws.Range("D39").Style = "Header4"
$ws.Range("E39").Style = "Header4"

# ============================================================
# Section: sys_sec_group_owner_approve_user_membership_MS (rows 40-44)
# ============================================================

# Row 40: section title + short description
$ws.Range("A40").Value = "sys_sec_group_owner_approve_user_membership_MS"
$ws.Range("A40").Style = "Header1"
$ws.Range("B40").Style = "Header3"
$ws.Range("C40").Value = "群主批准用户入群"
$ws.Range("C40").Style = "Header2"
$ws.Rows.Item(40).RowHeight = 21

# Row 41: "Parameters" label row
$ws.Range("A41").Value = "Parameters"
$ws.Range("A41").Style = "Bold12"
$ws.Range("B41").Style = "Header3"
$ws.Range("C41").Style = "Header3"
$ws.Range("D41").Style = "Header4"
$ws.Range("E41").Style = "Header4"

# Row 42: Parameter 1 - @groupOwnerUsertToken / 群主Token
$ws.Range("A42").Value = "    Parameter 1"
$ws.Range("A42").Style = "Header3"
$ws.Range("B42").Value = "@groupOwnerUsertToken"
$ws.Range("B42").Style = "Header3"
$ws.Range("C42").Value = "群主Token"
$ws.Range("C42").Style = "Header3"
$ws.Range("D42").Style = "Header4"
$ws.Range("E42").Style = "Header4"

# Row 43: Parameter 2 - @applicationUsertToken / 申请用户Token
$ws.Range("A43").Value = "    Parameter 2"
$ws.Range("A43").Style = "Header3"
$ws.Range("B43").Value = "@applicationUsertToken"
$ws.Range("C43").Value = "申请用户Token"
$ws.Range("C43").Style = "Header3"

# Row 44: Parameter 3 - @groupName / 群名字
$ws.Range("A44").Value = "    Parameter 3"
$ws.Range("A44").Style = "Header3"
$ws.Range("B44").Value = "@groupName"
$ws.Range("C44").Value = "群名字"

# --- update view: scroll position / selection, matching the diff ---
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("E41").Select()
